$d = $word.ActiveDocument
$firstPara = $d.Paragraphs.First
$r = $firstPara.Range
$r.InsertBefore("Greet participant and show to the meeting room`r")
